# Apply crypto price/volume updates as described in the commit diff.
# Values are stored as text in the workbook (inline/shared strings), so we
# force text formatting before assignment to stop Excel from silently
# re-interpreting numeric-looking strings (e.g. "575.22") as real numbers,
# then restore the default "Normal" style so no stray number format sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '62.844.26'
Set-TextValue "E2" '  +0.37%  '
Set-TextValue "D3" '2.463.92'
Set-TextValue "E3" '  +0.73%  '
Set-TextValue "D5" '575.22'
Set-TextValue "E5" '  -0.62%  '
Set-TextValue "D6" '145.97'
Set-TextValue "E6" '  +0.37%  '
Set-TextValue "E8" '  -0.10%  '
Set-TextValue "D9" '2.463.02'
Set-TextValue "E9" '  +0.76%  '
Set-TextValue "E10" '  +1.32%  '
Set-TextValue "E11" '  +1.50%  '
Set-TextValue "E12" '  +0.89%  '
Set-TextValue "D13" '0.355'
Set-TextValue "E13" '  +0.68%  '
Set-TextValue "D14" '28.97'
Set-TextValue "E14" '  +1.92%  '
Set-TextValue "E15" '  -0.05%  '
Set-TextValue "D16" '2.910.00'
Set-TextValue "D17" '62.809.60'
Set-TextValue "E17" '  +0.33%  '
Set-TextValue "D18" '2.465.79'
Set-TextValue "E18" '  +1.09%  '
Set-TextValue "E19" '  +2.41%  '
Set-TextValue "D20" '11.02'
Set-TextValue "E20" '  +0.90%  '
Set-TextValue "D21" '327.33'
Set-TextValue "E21" '  +0.42%  '
Set-TextValue "D22" '2.24'
Set-TextValue "E22" '  +10.57%  '
Set-TextValue "E24" '  -0.01%  '
Set-TextValue "D25" '10.22'
Set-TextValue "E25" '  +19.48%  '
Set-TextValue "D26" '65.88'
Set-TextValue "E26" '  +0.76%  '
Set-TextValue "D27" '653.61'
Set-TextValue "E27" '  +1.28%  '
Set-TextValue "D28" '0.0₃0981'
Set-TextValue "E28" '  +0.39%  '
Set-TextValue "D29" '2.583.88'
Set-TextValue "E29" '  +0.85%  '
Set-TextValue "D30" '0.997'
Set-TextValue "E30" '  -13.78%  '
Set-TextValue "E31" '  +2.62%  '
Set-TextValue "D32" '8.00'
Set-TextValue "E32" '  -2.27%  '
Set-TextValue "E33" '  -0.80%  '
Set-TextValue "E34" '  -3.98%  '
Set-TextValue "E36" '  +2.96%  '
Set-TextValue "D37" '4.75'
Set-TextValue "E37" '  +0.40%  '
Set-TextValue "D38" '0.369'
Set-TextValue "E38" '  -0.98%  '
Set-TextValue "D39" '18.72'
Set-TextValue "E39" '  +0.74%  '
Set-TextValue "D40" '5.39'
Set-TextValue "E40" '  -1.30%  '
Set-TextValue "D41" '151.02'
Set-TextValue "E41" '  -1.79%  '
Set-TextValue "D42" '2.75'
Set-TextValue "E42" '  +1.99%  '
Set-TextValue "E43" '  -1.22%  '
Set-TextValue "D44" '0.0₆0315'
Set-TextValue "E44" '  -81.22%  '
Set-TextValue "E45" '  -0.02%  '
Set-TextValue "D46" '154.05'
Set-TextValue "E46" '  +6.74%  '
Set-TextValue "D48" '3.58'
Set-TextValue "E48" '  -0.16%  '
Set-TextValue "D49" '20.31'
Set-TextValue "E49" '  -1.06%  '
Set-TextValue "E50" '  +0.59%  '
Set-TextValue "E51" '  +0.03%  '

Write-Host "Updated 74 cells"
